# Apply the dated-worksheet content update (date stamp + 25 division problems).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-04 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-05 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("67÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=", 2) | Out-Null
$d.Content.Find.Execute("88÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷5=", 2) | Out-Null
$d.Content.Find.Execute("56÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=", 2) | Out-Null
$d.Content.Find.Execute("41÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷2=", 2) | Out-Null
$d.Content.Find.Execute("63÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=", 2) | Out-Null
$d.Content.Find.Execute("66÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷6=", 2) | Out-Null
$d.Content.Find.Execute("12÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=", 2) | Out-Null
$d.Content.Find.Execute("28÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=", 2) | Out-Null
$d.Content.Find.Execute("37÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=", 2) | Out-Null
$d.Content.Find.Execute("31÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=", 2) | Out-Null
$d.Content.Find.Execute("40÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=", 2) | Out-Null
$d.Content.Find.Execute("47÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=", 2) | Out-Null
$d.Content.Find.Execute("92÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2) | Out-Null
$d.Content.Find.Execute("79÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=", 2) | Out-Null
$d.Content.Find.Execute("54÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=", 2) | Out-Null
$d.Content.Find.Execute("65÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷8=", 2) | Out-Null
$d.Content.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=", 2) | Out-Null
$d.Content.Find.Execute("78÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷3=", 2) | Out-Null
$d.Content.Find.Execute("40÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=", 2) | Out-Null
$d.Content.Find.Execute("73÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷6=", 2) | Out-Null
$d.Content.Find.Execute("63÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=", 2) | Out-Null
$d.Content.Find.Execute("59÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=", 2) | Out-Null
$d.Content.Find.Execute("97÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=", 2) | Out-Null
$d.Content.Find.Execute("88÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷7=", 2) | Out-Null
$d.Content.Find.Execute("91÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=", 2) | Out-Null
